# Update Bgn-Tlr2 ligand-receptor pair table to include the third condition
# ("sCs" cluster) following Dr Hou advice: data now spans rows 2-10 (ECs/FAPs/sCs x ECs/FAPs/sCs).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    ,@("ECs", "Bgn", "Tlr2", "ECs", 3, 1, 41.47959733333333, 124.438792, 0.03392314276466685, 0.03392314276466685, 3, 1, 57.24915866666667, 171.747476, 0.9704198736548433, 0.9704198736548435, 2374.672049165443, 21372.04844248899, 0.03291969191566322, 0.03291969191566322)
    ,@("ECs", "Bgn", "Tlr2", "FAPs", 3, 1, 41.47959733333333, 124.438792, 0.03392314276466685, 0.03392314276466685, 3, 1, 1.569166, 4.707498, 0.02659864191768634, 0.02659864191768634, 65.08837382915733, 585.7953644624159, 0.0009023095271199257, 0.0009023095271199256)
    ,@("ECs", "Bgn", "Tlr2", "sCs", 3, 1, 41.47959733333333, 124.438792, 0.03392314276466685, 0.03392314276466685, 2, 0.6666666666666666, 0.1758903333333333, 0.527671, 0.002981484427470275, 0.002981484427470276, 7.295860201492444, 65.66274181343199, 0.0001011413218837052, 0.0001011413218837052)
    ,@("FAPs", "Bgn", "Tlr2", "ECs", 3, 1, 1112.909261, 3338.727783, 0.9101674599595009, 0.9101674599595008, 3, 1, 57.24915866666667, 171.747476, 0.9704198736548433, 0.9704198736548435, 63713.11886459174, 573418.0697813256, 0.8832445914986484, 0.8832445914986485)
    ,@("FAPs", "Bgn", "Tlr2", "FAPs", 3, 1, 1112.909261, 3338.727783, 0.9101674599595009, 0.9101674599595008, 3, 1, 1.569166, 4.707498, 0.02659864191768634, 0.02659864191768634, 1746.339373446326, 15717.05436101693, 0.02420921835259288, 0.02420921835259288)
    ,@("FAPs", "Bgn", "Tlr2", "sCs", 3, 1, 1112.909261, 3338.727783, 0.9101674599595009, 0.9101674599595008, 2, 0.6666666666666666, 0.1758903333333333, 0.527671, 0.002981484427470275, 0.002981484427470276, 195.7499808870437, 1761.749827983393, 0.002713650108259427, 0.002713650108259427)
    ,@("sCs", "Bgn", "Tlr2", "ECs", 3, 1, 68.36333833333333, 205.090015, 0.05590939727583234, 0.05590939727583233, 3, 1, 57.24915866666667, 171.747476, 0.9704198736548433, 0.9704198736548435, 3913.743603228015, 35223.69242905214, 0.05425559024053166, 0.05425559024053166)
    ,@("sCs", "Bgn", "Tlr2", "FAPs", 3, 1, 68.36333833333333, 205.090015, 0.05590939727583234, 0.05590939727583233, 3, 1, 1.569166, 4.707498, 0.02659864191768634, 0.02659864191768634, 107.2734261591633, 965.46083543247, 0.001487114037973532, 0.001487114037973532)
    ,@("sCs", "Bgn", "Tlr2", "sCs", 3, 1, 68.36333833333333, 205.090015, 0.05590939727583234, 0.05590939727583233, 2, 0.6666666666666666, 0.1758903333333333, 0.527671, 0.002981484427470275, 0.002981484427470276, 12.02445036722944, 108.220053305065, 0.0001666929973271431, 0.0001666929973271431)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $rowData = $rows[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}
